$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.261.61"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "1.599.52"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("E8").Value = "  -0.50%  "

$ws.Range("E9").Value = "  +0.14%  "

$ws.Range("E10").Value = "  -1.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.92%  "

$ws.Range("D12").Value = "1.825.21"
$ws.Range("E12").Value = "  +0.70%  "

$ws.Range("D13").Value = "1.606.76"
$ws.Range("E13").Value = "  +1.25%  "

$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("E15").Value = "  -2.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.58%  "

$ws.Range("D17").Value = "26.278.16"
$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.70%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "

$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("E28").Value = "  +1.05%  "

$ws.Range("E29").Value = "  +1.76%  "

$ws.Range("E30").Value = "  -0.64%  "

$ws.Range("E31").Value = "  +0.61%  "

$ws.Range("E32").Value = "  -0.12%  "

$ws.Range("D33").Value = "1.444.94"
$ws.Range("E33").Value = "  +4.22%  "

$ws.Range("E34").Value = "  +0.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("E36").Value = "  +0.60%  "

$ws.Range("E37").Value = "  -2.52%  "

$ws.Range("E38").Value = "  -1.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.821"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.20%  "

$ws.Range("E40").Value = "  -1.29%  "

$ws.Range("E41").Value = "  +0.25%  "

$ws.Range("E42").Value = "  +2.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.924"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.26%  "

$ws.Range("D44").Value = "1.737.32"
$ws.Range("E44").Value = "  +0.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.760"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.32%  "

$ws.Range("E48").Value = "  +0.17%  "

$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("E50").Value = "  -3.19%  "

$ws.Range("E51").Value = "  +0.12%  "
